$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LocalData")

# Revert a test cutscene row: delete row 19 (CutsceneData.Cutscene.1000013.10),
# shifting all subsequent rows up by one.
$ws.Rows.Item(19).Delete()
